$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The existing single data row (row 2) is being replaced, and five new rows
# (3-7) are being added, expanding the Sending/Target cluster combinations
# (FAPs/sCs sending clusters x ECs/M2/sCs target clusters) per the new
# "Dr Hou advice" pairing logic. Clear the old row first so stale values
# don't linger on columns that a new row wouldn't otherwise touch.
$ws.Range("A2:T2").ClearContents()

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol,
# D=Target cluster, E..T = the numeric metric columns.
$newRows = @(
    @("FAPs", "Epo", "Epor", "ECs", 3, 1, 0.419201, 1.257603, 0.8650221929663464, 0.8650221929663463, 2, 0.6666666666666666, 1.054790666666667, 3.164372, 0.5652431679977079, 0.5652431679977079, 0.4421693022573333, 3.979523720316, 0.4889478847406222, 0.4889478847406222),
    @("FAPs", "Epo", "Epor", "M2", 3, 1, 0.419201, 1.257603, 0.8650221929663464, 0.8650221929663463, 2, 0.6666666666666666, 0.1924463333333334, 0.577339, 0.103128496070825, 0.103128496070825, 0.08067369537966668, 0.726063258417, 0.08920843782850629, 0.08920843782850629),
    @("FAPs", "Epo", "Epor", "sCs", 3, 1, 0.419201, 1.257603, 0.8650221929663464, 0.8650221929663463, 3, 1, 0.618846, 1.856538, 0.3316283359314672, 0.3316283359314672, 0.259420862046, 2.334787758414, 0.2868658703972179, 0.2868658703972179),
    @("sCs", "Epo", "Epor", "ECs", 1, 0.3333333333333333, 0.065412, 0.196236, 0.1349778070336537, 0.1349778070336536, 2, 0.6666666666666666, 1.054790666666667, 3.164372, 0.5652431679977079, 0.5652431679977079, 0.068995967088, 0.6209637037919999, 0.07629528325708571, 0.07629528325708569),
    @("sCs", "Epo", "Epor", "M2", 1, 0.3333333333333333, 0.065412, 0.196236, 0.1349778070336537, 0.1349778070336536, 2, 0.6666666666666666, 0.1924463333333334, 0.577339, 0.103128496070825, 0.103128496070825, 0.012588299556, 0.113294696004, 0.01392005824231873, 0.01392005824231873),
    @("sCs", "Epo", "Epor", "sCs", 1, 0.3333333333333333, 0.065412, 0.196236, 0.1349778070336537, 0.1349778070336536, 3, 1, 0.618846, 1.856538, 0.3316283359314672, 0.3316283359314672, 0.040479954552, 0.364319590968, 0.04476246553424925, 0.04476246553424924)
)

$rowCount = $newRows.Length
$colCount = 20

# Write column-by-column (all rows for column A, then all rows for column B, ...)
for ($c = 0; $c -lt $colCount; $c++) {
    for ($r = 0; $r -lt $rowCount; $r++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $newRows[$r][$c]
    }
}
